$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the correct MSRP value for the SS_TEST product row (G3): it should
# reference the same "MSRP_2000" value as the other row instead of the
# stale "MSRP_2000_SS_TEST" placeholder.
$ws.Range("G3").Value = "MSRP_2000"

# Update the last active cell/selection recorded in the sheet view.
$ws.Range("I10").Select()
